# Operation paused due to management decision.
#
# Adds 14 new "REF_PVLookup2_*_Arr" rows to the "IO Mapping" sheet
# (mirroring the existing "REF_PVLookup_*_Arr" block), and updates the
# active-sheet/selection view state so "IO Mapping" becomes the
# selected tab instead of "Pump".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IO Mapping")

# Column C / D text used by every row in this array-variable block.
$arrType    = "ARRAY [14] OF WORD"
$arrDefault = "[14(0)]"

# (row, A: variable name, B: PLC address, F: description)
$newRows = @(
    @(220, "REF_PVLookup2_1_Arr",  "D5120", "PV look-up values when only 1 valves are on"),
    @(221, "REF_PVLookup2_2_Arr",  "D5134", "PV look-up values when only 2 valves are on"),
    @(222, "REF_PVLookup2_3_Arr",  "D5148", "PV look-up values when only 3 valves are on"),
    @(223, "REF_PVLookup2_4_Arr",  "D5162", "PV look-up values when only 4 valves are on"),
    @(224, "REF_PVLookup2_5_Arr",  "D5176", "PV look-up values when only 5 valves are on"),
    @(225, "REF_PVLookup2_6_Arr",  "D5190", "PV look-up values when only 6 valves are on"),
    @(226, "REF_PVLookup2_7_Arr",  "D5204", "PV look-up values when only 7 valves are on"),
    @(227, "REF_PVLookup2_8_Arr",  "D5218", "PV look-up values when only 8 valves are on"),
    @(228, "REF_PVLookup2_9_Arr",  "D5232", "PV look-up values when only 9 valves are on"),
    @(229, "REF_PVLookup2_10_Arr", "D5246", "PV look-up values when only 10 valves are on"),
    @(230, "REF_PVLookup2_11_Arr", "D5260", "PV look-up values when only 11 valves are on"),
    @(231, "REF_PVLookup2_12_Arr", "D5274", "PV look-up values when only 12 valves are on"),
    @(232, "REF_PVLookup2_13_Arr", "D5288", "PV look-up values when only 13 valves are on"),
    @(233, "REF_PVLookup2_14_Arr", "D5302", "PV look-up values when only 14 valves are on")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]      # A: variable name
    $ws.Cells.Item($r, 2).Value = $row[2]      # B: PLC address

    $cCell = $ws.Cells.Item($r, 3)             # C: type
    $cCell.Value = $arrType
    $cCell.HorizontalAlignment = -4131         # xlHAlignLeft

    $dCell = $ws.Cells.Item($r, 4)             # D: default value
    $dCell.Value = $arrDefault
    $dCell.HorizontalAlignment = -4152         # xlHAlignRight

    $ws.Cells.Item($r, 6).Value = $row[3]      # F: description
}

# Make "IO Mapping" the active sheet/tab and restore the author's last
# selection + scroll position (was on "Pump" before the edit).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 198
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C236").Select()
